# Update column G ("K") values in the save_data sheet for sborz_josh.
# These values were regenerated to use K (strikeouts) instead of Strike# (pitch-strike count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 5
    7  = 0
    8  = 2
    9  = 2
    10 = 1
    11 = 3
    12 = 0
    13 = 2
    14 = 2
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 2
    22 = 3
    23 = 2
    24 = 2
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
